# Actualización automática 2025-11-24 16:30:09
# Updates sales figures for client "AGUILAR REYES CESAR VINICIO" (+172.99 in
# 240X80 PORCELANATO, flowing into PORCELANATO group total as +158.98),
# client "FAREZ ARCINIEGA DIANA FRANCISCA" (+12.24 in PANELES DECORATIVOS),
# and client "MATUTE GUANOLIQUE DOLORES MATILDE" (+25.2 in LAVABOS), across
# the three report sheets that mirror the same underlying figures.

$wb = $excel.ActiveWorkbook

$wsGrupo   = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumpl   = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Sheet "VENTAS POR GRUPO" ---------------------------------------------
# AGUILAR REYES CESAR VINICIO (row 4): 240X80 PORCELANATO (D) and PORCELANATO (M)
$wsGrupo.Range("D4").Value = 734.6900000000001
$wsGrupo.Range("M4").Value = 6677.89

# FAREZ ARCINIEGA DIANA FRANCISCA (row 22): PANELES DECORATIVOS (K)
$wsGrupo.Range("K22").Value = 113.76

# MATUTE GUANOLIQUE DOLORES MATILDE (row 31): LAVABOS (I)
$wsGrupo.Range("I31").Value = 104.4

# --- Sheet "VENTA MENSUAL" (column F = noviembre) --------------------------
$wsMensual.Range("F4").Value = 9211.93
$wsMensual.Range("F22").Value = 140.76
$wsMensual.Range("F31").Value = 104.4

# Row 60 is the column total row - updated to reflect the new November sum
$wsMensual.Range("F60").Value = 62078.81

# --- Sheet "CUMPLIMIENTO MENSUAL" ------------------------------------------
# Row 3: 240X80 PORCELANATO
$wsCumpl.Range("D3").Value = 4467.17
$wsCumpl.Range("E3").Value = 2156.09
$wsCumpl.Range("F3").Value = 0.6744669543397058

# Row 7: LAVABOS
$wsCumpl.Range("D7").Value = 468
$wsCumpl.Range("E7").Value = 852
$wsCumpl.Range("F7").Value = 0.3545454545454546

# Row 10: PANELES DECORATIVOS
$wsCumpl.Range("D10").Value = 823.89
$wsCumpl.Range("E10").Value = 3488.11
$wsCumpl.Range("F10").Value = 0.1910691094619666

# Row 12: PORCELANATO
$wsCumpl.Range("D12").Value = 34780.16
$wsCumpl.Range("E12").Value = 30163.84
$wsCumpl.Range("F12").Value = 0.5355407735895541

# Row 14: TOTAL
$wsCumpl.Range("D14").Value = 60204.9
$wsCumpl.Range("E14").Value = 38751.35685923838
$wsCumpl.Range("F14").Value = 0.6083991241265244

Write-Output "edit.ps1 applied"
